# Auto-generated script applying the Leviathan_Profits.xlsx numeric diff
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of the active workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 217.75
$ws.Range("I2").Value = 232
$ws.Range("K2").Value = 232
$ws.Range("M2").Value = -119
$ws.Range("H92").Value = 508.5
$ws.Range("I92").Value = 463.9091
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 463.9091
$ws.Range("L92").Value = 999
$ws.Range("M92").Value = 784.0908999999999
$ws.Range("N92").Value = -3495
$ws.Range("H100").Value = 5890.4
$ws.Range("I100").Value = 2304.5454
$ws.Range("K100").Value = 2304.5454
$ws.Range("M100").Value = -1763.5454
$ws.Range("H132").Value = 9002.777
$ws.Range("I132").Value = 2924.2856
$ws.Range("J132").Value = 30277.5
$ws.Range("K132").Value = 8772.856800000001
$ws.Range("L132").Value = 90832.5
$ws.Range("M132").Value = -6242.856800000001
$ws.Range("N132").Value = -95892.5
$ws.Range("H138").Value = 1640.8673
$ws.Range("I138").Value = 1230.591
$ws.Range("J138").Value = 1759.6316
$ws.Range("K138").Value = 3691.773
$ws.Range("L138").Value = 5278.8948
$ws.Range("M138").Value = 1448.227
$ws.Range("N138").Value = -15558.8948

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5098.4565
$ws.Range("I32").Value = 4525.2197
$ws.Range("K32").Value = 4525.2197
$ws.Range("M32").Value = -4238.2197
$ws.Range("H45").Value = 6601.407
$ws.Range("I45").Value = 9047.143
$ws.Range("K45").Value = 9047.143
$ws.Range("M45").Value = -8670.143
$ws.Range("H61").Value = 2769.0667
$ws.Range("I61").Value = 1329
$ws.Range("K61").Value = 1329
$ws.Range("M61").Value = -1117
$ws.Range("H74").Value = 2650.8
$ws.Range("I74").Value = 2480.5
$ws.Range("K74").Value = 2480.5
$ws.Range("M74").Value = -1606.5
$ws.Range("H77").Value = 2650.8
$ws.Range("I77").Value = 2480.5
$ws.Range("K77").Value = 12402.5
$ws.Range("M77").Value = -8034.5
$ws.Range("H102").Value = 3208.3333
$ws.Range("I102").Value = 3208.3333
$ws.Range("K102").Value = 3208.3333
$ws.Range("M102").Value = -1586.3333
$ws.Range("H132").Value = 2450.4443
$ws.Range("I132").Value = 1989.674
$ws.Range("K132").Value = 5969.022
$ws.Range("M132").Value = -3439.022
$ws.Range("H136").Value = 2769.0667
$ws.Range("I136").Value = 1329
$ws.Range("K136").Value = 3987
$ws.Range("M136").Value = -1437

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 114370
$ws.Range("I107").Value = 200366.4
$ws.Range("K107").Value = 200366.4
$ws.Range("M107").Value = -198446.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4453.9443
$ws.Range("I31").Value = 2319.7693
$ws.Range("K31").Value = 2319.7693
$ws.Range("M31").Value = -2024.7693
$ws.Range("H34").Value = 4453.9443
$ws.Range("I34").Value = 2319.7693
$ws.Range("K34").Value = 2319.7693
$ws.Range("M34").Value = -2117.7693
$ws.Range("H141").Value = 364554.88
$ws.Range("J141").Value = 442284.84
$ws.Range("L141").Value = 442284.84
$ws.Range("N141").Value = -452644.84

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 373.52173
$ws.Range("J2").Value = 167.54546
$ws.Range("L2").Value = 1005.27276
$ws.Range("N2").Value = -1231.27276
$ws.Range("H8").Value = 990.0833
$ws.Range("I8").Value = 990.0833
$ws.Range("K8").Value = 2970.2499
$ws.Range("M8").Value = -2831.2499
$ws.Range("H12").Value = 151.61539
$ws.Range("J12").Value = 178.18182
$ws.Range("L12").Value = 534.5454599999999
$ws.Range("N12").Value = -880.5454599999999
$ws.Range("H62").Value = 11499.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 11499.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 34498.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -35870.5
$ws.Range("H64").Value = 3777.625
$ws.Range("J64").Value = 4114.9
$ws.Range("L64").Value = 12344.7
$ws.Range("N64").Value = -12884.7
$ws.Range("H65").Value = 11499.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 11499.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 103495.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -110359.5
$ws.Range("H67").Value = 3777.625
$ws.Range("J67").Value = 4114.9
$ws.Range("L67").Value = 12344.7
$ws.Range("N67").Value = -14216.7
$ws.Range("H88").Value = 12497
$ws.Range("J88").Value = 12497
$ws.Range("L88").Value = 37491
$ws.Range("N88").Value = -38347
$ws.Range("H91").Value = 12497
$ws.Range("J91").Value = 12497
$ws.Range("L91").Value = 37491
$ws.Range("N91").Value = -40455
$ws.Range("H101").Value = 4000
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H105").Value = 14999
$ws.Range("J105").Value = 14999
$ws.Range("L105").Value = 44997
$ws.Range("N105").Value = -50239
$ws.Range("H107").Value = 662.2
$ws.Range("I107").Value = 449.8
$ws.Range("J107").Value = 692.54285
$ws.Range("K107").Value = 1349.4
$ws.Range("L107").Value = 2077.62855
$ws.Range("M107").Value = 570.5999999999999
$ws.Range("N107").Value = -5917.62855

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 47500.824
$ws.Range("I97").Value = 65979.336
$ws.Range("K97").Value = 65979.336
$ws.Range("M97").Value = -65483.336
$ws.Range("H126").Value = 2324.5386
$ws.Range("I126").Value = 2320.4546
$ws.Range("K126").Value = 6961.3638
$ws.Range("M126").Value = -4491.3638

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3544.0557
$ws.Range("I22").Value = 3591.2354
$ws.Range("J22").Value = 2742
$ws.Range("K22").Value = 3591.2354
$ws.Range("L22").Value = 2742
$ws.Range("M22").Value = -3296.2354
$ws.Range("N22").Value = -3332
$ws.Range("H27").Value = 3544.0557
$ws.Range("I27").Value = 3591.2354
$ws.Range("J27").Value = 2742
$ws.Range("K27").Value = 3591.2354
$ws.Range("L27").Value = 2742
$ws.Range("M27").Value = -3484.2354
$ws.Range("N27").Value = -2956
$ws.Range("H40").Value = 7077.108
$ws.Range("I40").Value = 6289.9165
$ws.Range("K40").Value = 6289.9165
$ws.Range("M40").Value = -6153.9165
$ws.Range("H93").Value = 15401.36
$ws.Range("I93").Value = 1908.8
$ws.Range("K93").Value = 1908.8
$ws.Range("M93").Value = -660.8
$ws.Range("H122").Value = 204865
$ws.Range("I122").Value = 669668
$ws.Range("K122").Value = 2009004
$ws.Range("M122").Value = -2006554

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 20006906
$ws.Range("I107").Value = 9766.308000000001
$ws.Range("J107").Value = 41670470
$ws.Range("K107").Value = 29298.924
$ws.Range("L107").Value = 125011410
$ws.Range("M107").Value = -27378.924
$ws.Range("N107").Value = -125015250
$ws.Range("H122").Value = 1692.9166
$ws.Range("I122").Value = 1756.8182
$ws.Range("K122").Value = 5270.4546
$ws.Range("M122").Value = -2820.4546
